$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("211:211").Insert()

$ws.Range("A211").Value = 5
$ws.Range("B211").Value = "Macroferia Regional de Talca"
$ws.Range("C211").Value = "Maule"
$ws.Range("D211").Value = 44719
$ws.Range("E211").Value = 7
$ws.Range("F211").Value = 100114014
$ws.Range("G211").Value = "Betarraga"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 2500
$ws.Range("K211").Value = 650
$ws.Range("L211").Value = 650
$ws.Range("M211").Value = 650
$ws.Range("N211").Value = "$/paquete 5 unidades"
$ws.Range("O211").Value = "Región del Maule"
$ws.Range("P211").Value = 130
$ws.Range("Q211").Value = 5
$ws.Range("R211").Value = "Hortaliza"
